# [Excel] (ExcelApi 1.14) Add worksheet protection change events code sample
#
# The "Snippets" table gains 5 new rows describing the new
# excel-events-worksheet-protection code sample:
#   - two rows for Class "Worksheet" (onProtectionChanged) inserted in
#     alphabetical order right before the existing "onRowSorted" row
#   - three rows for the new Class "WorksheetProtectionChangedEventArgs"
#     inserted in alphabetical order right before the existing
#     "WorksheetSingleClickedEventArgs" row
#
# The table auto-filter/used range/dimension grow from A1:E290 to A1:E295,
# and the view is left scrolled/selected near the newly-added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------
# 1) Insert the two new "Worksheet" / "onProtectionChanged" rows just
#    above the current row 266 ("Worksheet" / "onRowSorted"). Inserting
#    here copies the formatting of row 265 (the row above), which is the
#    same formatting the surrounding rows already use.
# ---------------------------------------------------------------------
$ws.Range("A266:A267").EntireRow.Insert()
$lo.Resize($ws.Range("A1:E292"))

$ws.Range("A266").Value2 = "Worksheet"
$ws.Range("B266").Value2 = "onProtectionChanged"
$ws.Range("D266").Value2 = "excel-events-worksheet-protection"
$ws.Range("E266").Value2 = "registerEvent"

$ws.Range("A267").Value2 = "Worksheet"
$ws.Range("B267").Value2 = "onProtectionChanged"
$ws.Range("D267").Value2 = "excel-events-worksheet-protection"
$ws.Range("E267").Value2 = "checkProtection"

# ---------------------------------------------------------------------
# 2) Insert the three new "WorksheetProtectionChangedEventArgs" rows
#    just above the row that now holds "WorksheetSingleClickedEventArgs"
#    (originally row 290, now shifted to row 292 by step 1).
#
#    Inserting new rows *above* that row would copy formatting from the
#    unstyled row above it, so instead new rows are appended *below* it
#    (inheriting its formatting) and then all four rows' values are
#    rewritten in the correct final order.
# ---------------------------------------------------------------------
$ws.Rows.Item(293).Insert()
$ws.Rows.Item(294).Insert()
$ws.Rows.Item(295).Insert()
$lo.Resize($ws.Range("A1:E295"))

$ws.Range("A292").Value2 = "WorksheetProtectionChangedEventArgs"
$ws.Range("B292").Value2 = "isProtected"
$ws.Range("C292").Value2 = ""
$ws.Range("D292").Value2 = "excel-events-worksheet-protection"
$ws.Range("E292").Value2 = "checkProtection"

$ws.Range("A293").Value2 = "WorksheetProtectionChangedEventArgs"
$ws.Range("B293").Value2 = "source"
$ws.Range("D293").Value2 = "excel-events-worksheet-protection"
$ws.Range("E293").Value2 = "checkProtection"

$ws.Range("A294").Value2 = "WorksheetProtectionChangedEventArgs"
$ws.Range("B294").Value2 = "worksheetId"
$ws.Range("D294").Value2 = "excel-events-worksheet-protection"
$ws.Range("E294").Value2 = "checkProtection"

$ws.Range("A295").Value2 = "WorksheetSingleClickedEventArgs"
$ws.Range("B295").Value2 = "address"
$ws.Range("D295").Value2 = "excel-event-worksheet-single-click"
$ws.Range("E295").Value2 = "registerClickHandler"

# ---------------------------------------------------------------------
# 3) Leave the view scrolled/selected near the new rows, matching the
#    final saved workbook (frozen header row still in place).
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
$win.ScrollRow = 274
$win.ScrollColumn = 1
$ws.Range("A282").Select()
